# Test - user can register ONLY with valid e-mail
# Rename sheets, add a new "TC Reg2" registration test case on Sheet3
# (renamed to TCReg2), and reference it from the Test Plan sheet.

$wb = $excel.ActiveWorkbook

$wsPlan = $wb.Worksheets.Item(1)   # "Test Plan"
$wsReg1 = $wb.Worksheets.Item(2)   # "TC Reg1"
$wsReg2 = $wb.Worksheets.Item(3)   # "Sheet3"

# ---------------------------------------------------------------
# 1. Rename the sheets
# ---------------------------------------------------------------
$wsPlan.Name = "TestPlan"
$wsReg1.Name = "TCReg1"
$wsReg2.Name = "TCReg2"

# ---------------------------------------------------------------
# 2. Build out the new "TCReg2" sheet content (mirrors TCReg1,
#    but for the "register only with valid e-mail" test case)
# ---------------------------------------------------------------

# Column widths, matching TCReg1
$wsReg2.Columns.Item(1).ColumnWidth = 16
$wsReg2.Columns.Item(2).ColumnWidth = 101
$wsReg2.Columns.Item(3).ColumnWidth = 38.5703125
$wsReg2.Columns.Item(4).ColumnWidth = 41.42578125
$wsReg2.Columns.Item(5).ColumnWidth = 18.28515625
$wsReg2.Columns.Item(6).ColumnWidth = 17.85546875

# Copy the header block / layout (rows 1-8) from TCReg1 so styles match.
$wsReg1.Range("A1:F8").Copy()
$wsReg2.Range("A1:F8").PasteSpecial(-4104)

# Row 1-2: identify the test case
$wsReg2.Range("B1").Value = "TC Reg2"
$wsReg2.Range("B2").Value = "Verify that user can register using only valid email adress"

# Row 9: step 2 - invalid e-mail
$wsReg2.Range("A9").Value = 2
$wsReg2.Range("B9").Value = "Enter invalid e-mail address"
$wsReg2.Range("C9").Value = "pera#gmail.com"
$wsReg2.Range("D9").Value = "E-mail is entered and visible"

# Row 10: step 3 - submit invalid e-mail
$wsReg2.Range("A10").Value = 3
$wsReg2.Range("B10").Value = "Click ""Create an account button"""
$wsReg2.Range("D10").Value = "User can`t registrer with invalid e-mail. Error message: ""Invalid email address."" is displayed."

# Row 11: step 4 - refresh page
$wsReg2.Range("A11").Value = 4
$wsReg2.Range("B11").Value = "Refresh page"
$wsReg2.Range("D11").Value = "Page is refreshed, error message is dissapeared. E-mail field is empty and ready for entering an e-mail."

# Row 12: step 5 - valid e-mail (with hyperlink on the test data cell)
$wsReg2.Range("A12").Value = 5
$wsReg2.Range("C12").Value = "pera@gmail.com"
$wsReg2.Range("B12").Value = "Enter valid e-mail address (for next text enter e-mail in format pera+n@gmail.com…n = 1, 2, 3...)"
$wsReg2.Range("D12").Value = "E-mail is entered and visible"

# Row 13: step 6 - submit valid e-mail
$wsReg2.Range("A13").Value = 6
$wsReg2.Range("B13").Value = "Click ""Create an account button"""
$wsReg2.Range("D13").Value = "E-mail is accepted, and user is forwarded to ""CREATE AN ACCOUNT"""

# Row heights for the wrapped, multi-line steps
$wsReg2.Rows.Item(10).RowHeight = 45
$wsReg2.Rows.Item(11).RowHeight = 45
$wsReg2.Rows.Item(13).RowHeight = 30

# Hyperlink for the valid e-mail test-data cell
$wsReg2.Hyperlinks.Add($wsReg2.Range("C12"), "mailto:pera@gmail.com")

# ---------------------------------------------------------------
# 3. Update the "TestPlan" sheet - rename the existing TC Reg1
#    reference and add a new row pointing at TC Reg2
# ---------------------------------------------------------------
$wsPlan.Range("C2").Value = "TCReg1"

$wsPlan.Range("A3").Value = ""
$wsPlan.Range("B3").Value = "TSu1"
$wsPlan.Range("C3").Value = "TCReg2"
$wsPlan.Range("D3").Value = "Verify that user can register using only valid email adress"

# ---------------------------------------------------------------
# 4. Restore selections / active sheet
# ---------------------------------------------------------------
$wsPlan.Select()
$wsPlan.Range("D13").Select()

$wsReg1.Select()
$wsReg1.Range("B36").Select()

$wsReg2.Select()
$wsReg2.Range("D13").Select()
